$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = newValue; E = newValue }  (only cells that changed)
$updates = @{
    2  = @{ D = "28.110.80";     E = "  +2.12%  " }
    3  = @{ D = "1.910.75";      E = "  +2.16%  " }
    5  = @{ D = "316.94";        E = "  +1.32%  " }
    6  = @{               E = "  -0.90%  " }
    7  = @{               E = "  +0.78%  " }
    8  = @{               E = "  +1.24%  " }
    9  = @{ D = "0.07363";       E = "  +0.07%  " }
    10 = @{ D = "0.9350";        E = "  -0.27%  " }
    11 = @{ D = "20.86";         E = "  +0.56%  " }
    12 = @{ D = "0.07793";       E = "  -0.69%  " }
    13 = @{ D = "1.910.94";      E = "  +2.34%  " }
    14 = @{ D = "5.522";         E = "  +1.48%  " }
    15 = @{ D = "6.653";         E = "  +0.99%  " }
    16 = @{ D = "91.72";         E = "  +1.02%  " }
    17 = @{               E = "  -0.89%  " }
    18 = @{ D = "0.000008836";   E = "  -0.86%  " }
    20 = @{ D = "28.147.94";     E = "  +2.25%  " }
    21 = @{               E = "  -0.62%  " }
    22 = @{ D = "5.173";         E = "  +0.62%  " }
    23 = @{ D = "2.166.80";      E = "  +2.10%  " }
    24 = @{ D = "10.93";         E = "  +1.66%  " }
    25 = @{ D = "156.13";        E = "  +1.24%  " }
    26 = @{ D = "1.922";         E = "  -2.12%  " }
    27 = @{ D = "18.54";         E = "  -0.04%  " }
    28 = @{ D = "2.118";         E = "  +4.66%  " }
    29 = @{ D = "116.68";        E = "  +0.61%  " }
    30 = @{ D = "4.964";         E = "  -0.73%  " }
    31 = @{ D = "0.08954";       E = "  +0.25%  " }
    32 = @{ D = "3.337";         E = "  -0.03%  " }
    33 = @{ D = "1.260";         E = "  +3.20%  " }
    34 = @{ D = "0.7771";        E = "  +2.69%  " }
    35 = @{ D = "4.690";         E = "  +1.68%  " }
    36 = @{ D = "2.650";         E = "  -1.45%  " }
    37 = @{ D = "0.02064";       E = "  +0.58%  " }
    38 = @{ D = "1.110";         E = "  -0.87%  " }
    39 = @{ D = "0.05327" }
    40 = @{ D = "0.5507";        E = "  +2.89%  " }
    41 = @{ D = "2.986";         E = "  -0.48%  " }
    42 = @{ D = "7.024";         E = "  -0.84%  " }
    43 = @{               E = "  +0.23%  " }
    44 = @{ D = "8.503";         E = "  +0.55%  " }
    45 = @{ D = "10.75";         E = "  +0.84%  " }
    46 = @{ D = "0.4849";        E = "  +0.59%  " }
    47 = @{ D = "108.34";        E = "  +5.40%  " }
    48 = @{               E = "  -0.97%  " }
    49 = @{ D = "1.660";         E = "  -0.07%  " }
    50 = @{ D = "68.31";         E = "  +1.30%  " }
    51 = @{               E = "  -0.28%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $cols["D"]
    }
    if ($cols.ContainsKey("E")) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $cols["E"]
    }
}
